$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2022" column (R), mirroring the formatting
# of the existing "2021" column (Q) for each of the three data rows.
$ws.Range("Q4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4122) | Out-Null
$ws.Range("R4").Value = 2022

$ws.Range("Q5").Copy() | Out-Null
$ws.Range("R5").PasteSpecial(-4122) | Out-Null
$ws.Range("R5").Value = 8.6821914120339212

$ws.Range("Q6").Copy() | Out-Null
$ws.Range("R6").PasteSpecial(-4122) | Out-Null
$ws.Range("R6").Value = 12.221423436376707

$excel.CutCopyMode = 0

# Match the authored view state: active cell parked one column past the
# new data (S4).
$ws.Range("S4").Select() | Out-Null
